$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like text cells (Y/AA) keep their literal text value instead of
# being auto-converted into a date serial number by Excel.
$dateTextCells = @("AA3", "Y3", "AA4", "Y4", "AA5", "Y5", "AA6", "Y6", "AA7", "Y7", "AA8", "Y8", "AA9", "Y9", "AA10", "Y10", "AA11", "Y11")
foreach ($addr in $dateTextCells) { $ws.Range($addr).NumberFormat = "@" }

# Row 3 <- content of original row 4
$ws.Range("A3").Value = 111780627
$ws.Range("B3").Value = 78604
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 6461
$ws.Range("F3").Value = "Norrlandslav"
$ws.Range("G3").Value = "Nephroma arcticum"
$ws.Range("H3").Value = "(L.) Torss."
$ws.Range("Q3").Value = 707647.2196405758
$ws.Range("R3").Value = 7397286.731778639
$ws.Range("Y3").Value = "2023-08-29"
$ws.Range("AA3").Value = "2023-08-29"
$ws.Range("AF3").ClearContents()

# Row 4 <- content of original row 5
$ws.Range("A4").Value = 111780621
$ws.Range("B4").Value = 56543
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 103021
$ws.Range("F4").Value = "Talltita"
$ws.Range("G4").Value = "Poecile montanus"
$ws.Range("H4").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q4").Value = 707631.1509720345
$ws.Range("R4").Value = 7397277.54798521
$ws.Range("Y4").Value = "2023-08-29"
$ws.Range("AA4").Value = "2023-08-29"
$ws.Range("AF4").ClearContents()

# Row 5 <- content of original row 3
$ws.Range("A5").Value = 111780624
$ws.Range("B5").Value = 95532
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 221945
$ws.Range("F5").Value = "Revlummer"
$ws.Range("G5").Value = "Lycopodium annotinum"
$ws.Range("H5").Value = "L."
$ws.Range("Q5").Value = 707600.9335272597
$ws.Range("R5").Value = 7397313.141869167
$ws.Range("Y5").Value = "2023-08-29"
$ws.Range("AA5").Value = "2023-08-29"
$ws.Range("AF5").ClearContents()

# Row 6 <- content of original row 7
$ws.Range("A6").Value = 111816142
$ws.Range("B6").Value = 78604
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 6461
$ws.Range("F6").Value = "Norrlandslav"
$ws.Range("G6").Value = "Nephroma arcticum"
$ws.Range("H6").Value = "(L.) Torss."
$ws.Range("Q6").Value = 707613.3456041727
$ws.Range("R6").Value = 7397270.22663033
$ws.Range("Y6").Value = "2023-08-29"
$ws.Range("AA6").Value = "2023-08-29"
$ws.Range("AF6").Value = ""

# Row 7 <- content of original row 6
$ws.Range("A7").Value = 111816118
$ws.Range("B7").Value = 78107
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6453
$ws.Range("F7").Value = "Vedskivlav"
$ws.Range("G7").Value = "Hertelidea botryosa"
$ws.Range("H7").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q7").Value = 707670.4513803272
$ws.Range("R7").Value = 7397327.948038339
$ws.Range("Y7").Value = "2023-08-22"
$ws.Range("AA7").Value = "2023-08-22"
$ws.Range("AF7").ClearContents()

# Row 8 <- content of original row 11
$ws.Range("A8").Value = 111816132
$ws.Range("B8").Value = 95532
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 221945
$ws.Range("F8").Value = "Revlummer"
$ws.Range("G8").Value = "Lycopodium annotinum"
$ws.Range("H8").Value = "L."
$ws.Range("Q8").Value = 707589.6730983062
$ws.Range("R8").Value = 7397240.139162621
$ws.Range("Y8").Value = "2023-08-22"
$ws.Range("AA8").Value = "2023-08-22"
$ws.Range("AF8").ClearContents()

# Row 9 <- content of original row 8
$ws.Range("A9").Value = 111816119
$ws.Range("B9").Value = 56543
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 103021
$ws.Range("F9").Value = "Talltita"
$ws.Range("G9").Value = "Poecile montanus"
$ws.Range("H9").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q9").Value = 707595.5401507822
$ws.Range("R9").Value = 7397262.905378895
$ws.Range("Y9").Value = "2023-08-22"
$ws.Range("AA9").Value = "2023-08-22"
$ws.Range("AF9").ClearContents()

# Row 10 <- content of original row 9
$ws.Range("A10").Value = 111816145
$ws.Range("B10").Value = 77597
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 864
$ws.Range("F10").Value = "Knottrig blåslav"
$ws.Range("G10").Value = "Hypogymnia bitteri"
$ws.Range("H10").Value = "(Lynge) Ahti"
$ws.Range("Q10").Value = 707626.9948496711
$ws.Range("R10").Value = 7397311.517900761
$ws.Range("Y10").Value = "2023-08-22"
$ws.Range("AA10").Value = "2023-08-22"
$ws.Range("AF10").Value = ""

# Row 11 <- content of original row 10
$ws.Range("A11").Value = 111816137
$ws.Range("B11").Value = 90658
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 4361
$ws.Range("F11").Value = "Orange taggsvamp"
$ws.Range("G11").Value = "Hydnellum aurantiacum"
$ws.Range("H11").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q11").Value = 707609.3988008115
$ws.Range("R11").Value = 7397264.348220735
$ws.Range("Y11").Value = "2023-08-22"
$ws.Range("AA11").Value = "2023-08-22"
$ws.Range("AF11").ClearContents()

# Restore default (General) style on the date-like text cells now that their
# literal text values are set, so no stray formatting is left behind.
foreach ($addr in $dateTextCells) { $ws.Range($addr).Style = "Normal" }
